$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row at 41. This shifts the old rows 41-46 down to 42-47 and
# Excel auto-adjusts the dependent formulas (F42/F43 refs move from F41/F40
# etc. to F42/F41 etc. automatically).
# ---------------------------------------------------------------------------
$ws.Rows.Item(41).Insert()
$ws.Rows.Item(41).RowHeight = 12.75

# ---------------------------------------------------------------------------
# Row 40: D40 picks up the same style as C40 (border/alignment tweak).
# ---------------------------------------------------------------------------
$ws.Cells.Item(40, 3).Copy($ws.Cells.Item(40, 4))
$ws.Cells.Item(40, 4).Value = 0

# ---------------------------------------------------------------------------
# Row 41 (brand-new item row, mirrors the other merit/demerit item rows
# such as row 40): A (blank, bordered), B (blank, bordered), C (blank text,
# bordered + readingOrder), D/E (0, bordered), F (blank, bordered), G (blank).
# ---------------------------------------------------------------------------
$ws.Cells.Item(40, 1).Copy($ws.Cells.Item(41, 1))
$ws.Cells.Item(40, 2).Copy($ws.Cells.Item(41, 2))
$ws.Cells.Item(41, 2).ClearContents()
$ws.Cells.Item(40, 3).Copy($ws.Cells.Item(41, 3))
$ws.Cells.Item(41, 3).Value = ""
$ws.Cells.Item(34, 4).Copy($ws.Cells.Item(41, 4))
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(42, 5).Copy($ws.Cells.Item(41, 5))
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(40, 6).Copy($ws.Cells.Item(41, 6))
$ws.Cells.Item(41, 6).ClearContents()
$ws.Cells.Item(40, 7).Copy($ws.Cells.Item(41, 7))

# ---------------------------------------------------------------------------
# Row 42 (was old row 41, the "JUMLAH" totals row): give it the JUMLAH label
# and extend the total formula to include the new row 41 merit/demerit pair.
# ---------------------------------------------------------------------------
$ws.Cells.Item(42, 2).Value = "JUMLAH"
$ws.Cells.Item(42, 6).Formula = "=D34-E34+D35-E35+D36-E36+D37-E37+D38-E38+D39-E39+D41-E41+D40-E40"

# ---------------------------------------------------------------------------
# Row 43 (was old row 42, the closing-balance row): formula already shifted
# automatically by the row insert (F20+F32+F26+F42); style indices bump by
# one step because of the row shift, so copy from an equivalent later row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(45, 2).Copy($ws.Cells.Item(43, 2))
$ws.Cells.Item(43, 2).Value = "CLOSING BALANCE FOR YEAR 2022 / BAKI PENUTUP TAHUN 2022"
$ws.Cells.Item(45, 6).Copy($ws.Cells.Item(43, 6))
$ws.Cells.Item(43, 6).Formula = "=F20+F32+F26+F42"

# ---------------------------------------------------------------------------
# Row 45/46 (was old rows 44/45): style bumps up one notch.
# ---------------------------------------------------------------------------
$ws.Cells.Item(46, 2).Copy($ws.Cells.Item(45, 2))
$ws.Cells.Item(45, 2).Value = "End Of 2022 Statement / Penyata 2022 Tamat "

$ws.Cells.Item(46, 2).Value = "May You Continue to Rise Higher in 2023 "

# ---------------------------------------------------------------------------
# Row 47 (was old row 46, the bottom border row): style bumps up one notch
# across the whole row.
# ---------------------------------------------------------------------------
for ($col = 1; $col -le 7; $col++) {
    $ws.Cells.Item(47, $col).Copy($ws.Cells.Item(47, $col))
}

Write-Host "done"
